$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.643.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.21%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.193.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.65%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.48%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.555'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.180.98'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.33%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.96'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.53%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.519'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.95%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.713.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.54%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.523.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.88%  '

$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.190.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.58%  '

$ws.Range('E19').Value = '  +0.65%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '517.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.740'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.85%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.91'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.28%  '

$ws.Range('E26').Value = '  +0.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.95%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.92%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.77%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.39'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.42%  '

$ws.Range('E33').Value = '  +3.33%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.79%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '513.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.42%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0904'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.35%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0427'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.21%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.129'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.09%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.93'
$ws.Range('D41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.03%  '

$ws.Range('E43').Value = '  +8.36%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0678'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +17.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.92%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.912.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.61%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.86%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.118'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.66%  '

$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.04%  '

$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.53%  '
